$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3428.2307
$ws.Range("I76").Value = 3380.5833
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 3380.5833
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -3065.5833
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 3428.2307
$ws.Range("I79").Value = 3380.5833
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 3380.5833
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -2288.5833
$ws.Range("N79").Value = -6184
$ws.Range("H86").Value = 3122.1562
$ws.Range("I86").Value = 4400.1665
$ws.Range("K86").Value = 4400.1665
$ws.Range("M86").Value = -3277.1665
$ws.Range("H89").Value = 3122.1562
$ws.Range("I89").Value = 4400.1665
$ws.Range("K89").Value = 22000.8325
$ws.Range("M89").Value = -16384.8325
$ws.Range("H107").Value = 544.9091
$ws.Range("I107").Value = 495.76923
$ws.Range("J107").Value = 615.8889
$ws.Range("K107").Value = 495.76923
$ws.Range("L107").Value = 615.8889
$ws.Range("M107").Value = 1424.23077
$ws.Range("N107").Value = -4455.8889
$ws.Range("H113").Value = 3168.5518
$ws.Range("I113").Value = 3187.7856
$ws.Range("J113").Value = 3150.6
$ws.Range("K113").Value = 3187.7856
$ws.Range("L113").Value = 3150.6
$ws.Range("M113").Value = 66.21439999999984
$ws.Range("N113").Value = -9658.6
$ws.Range("H129").Value = 530.0769
$ws.Range("I129").Value = 270
$ws.Range("J129").Value = 946.2
$ws.Range("K129").Value = 810
$ws.Range("L129").Value = 2838.6
$ws.Range("M129").Value = 4190
$ws.Range("N129").Value = -12838.6
$ws.Range("H132").Value = 1452.1842
$ws.Range("I132").Value = 1275.2162
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 3825.6486
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -1295.6486
$ws.Range("N132").Value = -29060
$ws.Range("H138").Value = 2254.4844
$ws.Range("I138").Value = 1251.7812
$ws.Range("J138").Value = 3257.1875
$ws.Range("K138").Value = 3755.3436
$ws.Range("L138").Value = 9771.5625
$ws.Range("M138").Value = 1384.6564
$ws.Range("N138").Value = -20051.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1766.5555
$ws.Range("I45").Value = 1674.875
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 1674.875
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -1297.875
$ws.Range("N45").Value = -3254
$ws.Range("H61").Value = 1901.2727
$ws.Range("I61").Value = 1434.9445
$ws.Range("J61").Value = 3999.75
$ws.Range("K61").Value = 1434.9445
$ws.Range("L61").Value = 3999.75
$ws.Range("M61").Value = -1222.9445
$ws.Range("N61").Value = -4423.75
$ws.Range("H110").Value = 1630
$ws.Range("I110").Value = 1200
$ws.Range("J110").Value = 1737.5
$ws.Range("K110").Value = 1200
$ws.Range("L110").Value = 1737.5
$ws.Range("M110").Value = 845
$ws.Range("N110").Value = -5827.5
$ws.Range("H132").Value = 2927
$ws.Range("I132").Value = 2437.8215
$ws.Range("J132").Value = 4448.8887
$ws.Range("K132").Value = 7313.4645
$ws.Range("L132").Value = 13346.6661
$ws.Range("M132").Value = -4783.4645
$ws.Range("N132").Value = -18406.6661
$ws.Range("H136").Value = 1901.2727
$ws.Range("I136").Value = 1434.9445
$ws.Range("J136").Value = 3999.75
$ws.Range("K136").Value = 4304.833500000001
$ws.Range("L136").Value = 11999.25
$ws.Range("M136").Value = -1754.833500000001
$ws.Range("N136").Value = -17099.25
$ws.Range("H138").Value = 73408.664
$ws.Range("J138").Value = 73408.664
$ws.Range("L138").Value = 73408.664
$ws.Range("N138").Value = -83688.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 870.3077
$ws.Range("I16").Value = 880.1111
$ws.Range("J16").Value = 848.25
$ws.Range("K16").Value = 880.1111
$ws.Range("L16").Value = 848.25
$ws.Range("M16").Value = -593.1111
$ws.Range("N16").Value = -1422.25
$ws.Range("H113").Value = 870.3077
$ws.Range("I113").Value = 880.1111
$ws.Range("J113").Value = 848.25
$ws.Range("K113").Value = 880.1111
$ws.Range("L113").Value = 848.25
$ws.Range("M113").Value = 1289.8889
$ws.Range("N113").Value = -5188.25
$ws.Range("H122").Value = 4984.615
$ws.Range("I122").Value = 6740
$ws.Range("J122").Value = 3887.5
$ws.Range("K122").Value = 20220
$ws.Range("L122").Value = 11662.5
$ws.Range("M122").Value = -17770
$ws.Range("N122").Value = -16562.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 5374.9
$ws.Range("I63").Value = 1805.75
$ws.Range("J63").Value = 7754.3335
$ws.Range("K63").Value = 5417.25
$ws.Range("L63").Value = 23263.0005
$ws.Range("M63").Value = -4668.25
$ws.Range("N63").Value = -24761.0005
$ws.Range("H66").Value = 5374.9
$ws.Range("I66").Value = 1805.75
$ws.Range("J66").Value = 7754.3335
$ws.Range("K66").Value = 16251.75
$ws.Range("L66").Value = 69789.0015
$ws.Range("M66").Value = -12507.75
$ws.Range("N66").Value = -77277.0015
$ws.Range("H120").Value = 19217.857
$ws.Range("I120").Value = 980
$ws.Range("J120").Value = 26513
$ws.Range("K120").Value = 2940
$ws.Range("L120").Value = 79539
$ws.Range("M120").Value = 1898
$ws.Range("N120").Value = -89215

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1972
$ws.Range("I31").Value = 1972
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1972
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -1680
$ws.Range("H37").Value = 1972
$ws.Range("I37").Value = 1972
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1972
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -1695
$ws.Range("H80").Value = 1733.3334
$ws.Range("I80").Value = 1533.3334
$ws.Range("J80").Value = 1933.3334
$ws.Range("K80").Value = 1533.3334
$ws.Range("L80").Value = 1933.3334
$ws.Range("M80").Value = -535.3334
$ws.Range("N80").Value = -3929.3334
$ws.Range("H83").Value = 1733.3334
$ws.Range("I83").Value = 1533.3334
$ws.Range("J83").Value = 1933.3334
$ws.Range("K83").Value = 7666.666999999999
$ws.Range("L83").Value = 9666.666999999999
$ws.Range("M83").Value = -2674.666999999999
$ws.Range("N83").Value = -19650.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 434769.44
$ws.Range("J46").Value = 558903.5600000001
$ws.Range("L46").Value = 558903.5600000001
$ws.Range("N46").Value = -559279.5600000001
$ws.Range("H61").Value = 2065.138
$ws.Range("I61").Value = 1022.63635
$ws.Range("J61").Value = 2702.2222
$ws.Range("K61").Value = 1022.63635
$ws.Range("L61").Value = 2702.2222
$ws.Range("M61").Value = -820.63635
$ws.Range("N61").Value = -3106.2222
$ws.Range("H68").Value = 2788.5435
$ws.Range("I68").Value = 1740
$ws.Range("J68").Value = 2976.7437
$ws.Range("K68").Value = 1740
$ws.Range("L68").Value = 2976.7437
$ws.Range("M68").Value = -991
$ws.Range("N68").Value = -4474.7437
$ws.Range("H71").Value = 2788.5435
$ws.Range("I71").Value = 1740
$ws.Range("J71").Value = 2976.7437
$ws.Range("K71").Value = 8700
$ws.Range("L71").Value = 14883.7185
$ws.Range("M71").Value = -4956
$ws.Range("N71").Value = -22371.7185
$ws.Range("H82").Value = 2941.1765
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 3142.8572
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 3142.8572
$ws.Range("M82").Value = -1639
$ws.Range("N82").Value = -3864.8572
$ws.Range("H85").Value = 2941.1765
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 3142.8572
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 3142.8572
$ws.Range("M85").Value = -752
$ws.Range("N85").Value = -5638.8572
$ws.Range("H113").Value = 2065.138
$ws.Range("I113").Value = 1022.63635
$ws.Range("J113").Value = 2702.2222
$ws.Range("K113").Value = 1022.63635
$ws.Range("L113").Value = 2702.2222
$ws.Range("M113").Value = 1147.36365
$ws.Range("N113").Value = -7042.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 910377.5600000001
$ws.Range("I122").Value = 1001115.3
$ws.Range("K122").Value = 3003345.9
$ws.Range("M122").Value = -3000895.9
